$d = $word.ActiveDocument

# Add a new movie title, "Se beber não case ", as its own paragraph right
# after the "Sorria 2 " entry.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Sorria 2 `r") {
        $p.Range.InsertParagraphAfter()
        $p.Next().Range.Text = "Se beber não case "
        break
    }
}
